$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking values that must stay as text (force Text format to avoid Excel auto-numeric conversion) ---
$textCells = @("D5","D7","D8","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D23","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "314.24"
$ws.Range("D7").Value = "0.3892"
$ws.Range("D8").Value = "0.4038"
$ws.Range("D11").Value = "52.96"
$ws.Range("D12").Value = "0.08751"
$ws.Range("D13").Value = "25.38"
$ws.Range("D14").Value = "7.524"
$ws.Range("D15").Value = "0.00001355"
$ws.Range("D16").Value = "7.955"
$ws.Range("D18").Value = "98.50"
$ws.Range("D19").Value = "0.07110"
$ws.Range("D20").Value = "19.96"
$ws.Range("D21").Value = "7.269"
$ws.Range("D23").Value = "14.26"
$ws.Range("D25").Value = "2.982"
$ws.Range("D27").Value = "22.77"
$ws.Range("D28").Value = "162.09"
$ws.Range("D29").Value = "8.755"
$ws.Range("D30").Value = "137.07"
$ws.Range("D31").Value = "5.226"
$ws.Range("D33").Value = "0.08839"
$ws.Range("D34").Value = "7.388"
$ws.Range("D35").Value = "1.034"
$ws.Range("D36").Value = "1.960"
$ws.Range("D37").Value = "0.2749"
$ws.Range("D38").Value = "0.02921"
$ws.Range("D39").Value = "10.77"
$ws.Range("D40").Value = "14.25"
$ws.Range("D41").Value = "0.09129"
$ws.Range("D42").Value = "0.7906"
$ws.Range("D43").Value = "1.461"
$ws.Range("D44").Value = "16.79"
$ws.Range("D45").Value = "0.7216"
$ws.Range("D46").Value = "2.592"
$ws.Range("D47").Value = "4.206"
$ws.Range("D48").Value = "1.002"
$ws.Range("D49").Value = "1.339"
$ws.Range("D50").Value = "137.89"
$ws.Range("D51").Value = "91.16"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Plain text / percentage values (safe to assign directly) ---
$ws.Range("D2").Value = "24.601.25"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.691.08"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  +7.20%  "
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "1.690.01"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "24.591.83"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  -7.52%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  +15.69%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").Value = "1.874.95"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E38").Value = "  +7.17%  "
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  +4.37%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  +0.54%  "
